# Update date formatting (dot -> dash) in column B and bump a handful of
# "want to go" counts in column F. The same edits are applied identically
# to the "展览" and "全部类型" sheets (both carry duplicate data).

$wb = $excel.ActiveWorkbook

# Row -> new date string (column B) for rows 2..17
$dates = [ordered]@{
    2  = "2024-03-02"
    3  = "2024-03-09"
    4  = "2024-03-09"
    5  = "2024-03-10"
    6  = "2024-03-16"
    7  = "2024-03-16"
    8  = "2024-03-23"
    9  = "2024-03-23"
    10 = "2024-03-23"
    11 = "2024-03-24"
    12 = "2024-03-30"
    13 = "2024-03-30"
    14 = "2024-03-31"
    15 = "2024-04-04"
    16 = "2024-04-13"
    17 = "2024-04-13"
}

# Row -> new numeric value for column F
$counts = [ordered]@{
    5  = 21
    7  = 2609
    9  = 1642
    12 = 535
    15 = 57
    17 = 7
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $dates.Keys) {
        $cell = $ws.Cells.Item($row, 2)
        # Force text storage so "2024-03-02" isn't reinterpreted as a date
        # serial number (it would otherwise parse just like typing it into
        # Excel would) - the source keeps these as plain strings.
        $cell.NumberFormat = "@"
        $cell.Value = $dates[$row]
    }

    foreach ($row in $counts.Keys) {
        $ws.Cells.Item($row, 6).Value = $counts[$row]
    }
}
